$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header for column C ("audioFalse" -> "currentPhase")
$ws.Range("C1").Value = "currentPhase"

# Update condition values in column C for the "none" condition (train2P2)
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
